{"js": "// Apply Dari (fa) text corrections to the CrisisText video scripts document.\n// Each edit below is a precise, whole-string search-and-replace against a\n// unique source sentence, matching the unified diff exactly.\n\nconst edits = [\n  {\n    from: \"\u062f\u0631\u0633 \u0627\u0645\u0631\u0648\u0632 \u062f\u0631\u0628\u0627\u0631\u0647 \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0632\u0645\u0627\u0646\u06cc \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f.\",\n    to: \"\u062f\u0631\u0633 \u0627\u0645\u0631\u0648\u0632 \u062f\u0631\u0628\u0627\u0631\u0647 \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0632\u0645\u0627\u0646\u06cc \u06a9\u0647 \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f.\"\n  },\n  {\n    from: \"\u0627\u06cc\u0646\u062c\u0627 \u0633\u0647 \u062a\u0648\u0635\u06cc\u0647 \u0628\u0631\u0627\u06cc \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0648\u0642\u062a\u06cc \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f: \",\n    to: \"\u0627\u06cc\u0646\u062c\u0627 \u0633\u0647 \u062a\u0648\u0635\u06cc\u0647 \u0628\u0631\u0627\u06cc \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0648\u0642\u062a\u06cc \u06a9\u0647 \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f: \"\n  },\n  {\n    from: \"\u0628\u0647 \u0632\u0628\u0627\u0646 \u0633\u0627\u062f\u0647\u060c \u0628\u0647 \u0622\u0646\u0647\u0627 \u0628\u06af\u0648\u06cc\u06cc\u062f \u06a9\u0647 \u0634\u062e\u0635 \u0648\u0641\u0627\u062a \u06a9\u0631\u062f\u0647 \u0648 \u062f\u06cc\u06af\u0631 \u0628\u0627\u0632 \u0646\u062e\u0648\u0627\u0647\u062f \u06af\u0634\u062a.\",\n    to: \"\u0628\u0647 \u0632\u0628\u0627\u0646 \u0633\u0627\u062f\u0647\u060c \u0628\u0647 \u0622\u0646\u0647\u0627 \u0628\u06af\u0648\u06cc\u06cc\u062f \u0634\u062e\u0635 \u06a9\u0647 \u0648\u0641\u0627\u062a \u06a9\u0631\u062f\u0647 \u0648 \u062f\u06cc\u06af\u0631 \u0628\u0627\u0632 \u0646\u062e\u0648\u0627\u0647\u062f \u06af\u0634\u062a.\"\n  },\n  {\n    from: \"\u0645\u062a\u0648\u062c\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a \u0622\u0646\u200c \u0647\u0627 \u0628\u0627\u0634\u06cc\u062f \u0648 \u0622\u0646\u0686\u0647 \u0631\u0627 \u06a9\u0647 \u0645\u0634\u0627\u0647\u062f\u0647 \u0645\u06cc\u200c\u06a9\u0646\u06cc\u062f \u0628\u0647 \u0632\u0628\u0627\u0646 \u0628\u06cc\u0627\u0648\u0631\u06cc\u062f \u062a\u0627 \u0627\u062d\u0633\u0627\u0633 \u06a9\u0646\u0646\u062f \u06a9\u0647 \u06af\u067e \u0647\u0627\u06cc \u0634\u0627\u062a \u0634\u0646\u06cc\u062f\u0647 \u0645\u06cc\u200c\u0634\u0648\u062f.\",\n    to: \"\u0645\u062a\u0648\u062c\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a \u0622\u0646\u200c \u0647\u0627 \u0628\u0627\u0634\u06cc\u062f \u0648 \u0622\u0646\u0686\u0647 \u0631\u0627 \u06a9\u0647 \u0645\u0634\u0627\u0647\u062f\u0647 \u0645\u06cc\u200c\u06a9\u0646\u06cc\u062f \u0628\u0647 \u0632\u0628\u0627\u0646 \u0628\u06cc\u0627\u0648\u0631\u06cc\u062f \u062a\u0627 \u0627\u062d\u0633\u0627\u0633 \u06a9\u0646\u0646\u062f \u06a9\u0647 \u06af\u067e \u0647\u0627\u06cc \u0634\u0627\u0646 \u0634\u0646\u06cc\u062f\u0647 \u0645\u06cc\u200c\u0634\u0648\u062f.\"\n  },\n  {\n    from: \"\u0645\u0645\u06a9\u0646 \u0627\u0633\u062a \u0644\u0627\u0632\u0645 \u0628\u0627\u0634\u062f \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u06a9\u0645\u06a9 \u06a9\u0646\u06cc\u062f \u062a\u0627 \u062f\u0631\u0628\u0627\u0631\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a\u0634 \u0635\u062d\u0628\u062a \u06a9\u0646\u0646\u062f \u0648 \u06cc\u0627 \u062a\u0648\u062c\u0647 \u0622\u0646\u0647\u0627 \u0631\u0627 \u0645\u062c\u062f\u062f\u0627 \u0645\u062a\u0645\u0631\u06a9\u0632 \u06a9\u0646\u06cc\u062f.\",\n    to: \"\u0645\u0645\u06a9\u0646 \u0627\u0633\u062a \u0644\u0627\u0632\u0645 \u0628\u0627\u0634\u062f \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u06a9\u0645\u06a9 \u06a9\u0646\u06cc\u062f \u062a\u0627 \u062f\u0631\u0628\u0627\u0631\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a \u0634\u0627\u0646 \u0635\u062d\u0628\u062a \u06a9\u0646\u0646\u062f \u0648 \u06cc\u0627 \u062a\u0648\u062c\u0647 \u0622\u0646\u0647\u0627 \u0631\u0627 \u0645\u062c\u062f\u062f\u0627 \u0645\u062a\u0645\u0631\u06a9\u0632 \u06a9\u0646\u06cc\u062f.\"\n  },\n  {\n    from: \"\u0634\u0634\u0645\u06cc\u0646 \u062a\u0648\u0635\u06cc\u0647 \u0622\u0631\u0627\u0645 \u06a9\u0631\u062f\u0646 \u0637\u0641\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a.\",\n    to: \"\u062a\u0648\u0635\u06cc\u0647 \u0634\u0634\u0645 \u0622\u0631\u0627\u0645 \u06a9\u0631\u062f\u0646 \u0637\u0641\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a.\"\n  }\n];\n\nfor (const edit of edits) {\n  const results = context.document.body.search(edit.from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + edit.from);\n  }\n\n  for (const range of results.items) {\n    range.insertText(edit.to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply Dari (fa) text corrections to the CrisisText video scripts document.\n# Each edit is a precise, whole-string search-and-replace against a unique\n# source sentence, matching the unified diff exactly. We replace the text of\n# the whole (single-run) paragraph via Range.Text so Word's OOXML writer\n# keeps the run's xml:space=\"preserve\" attribute (needed for the edits that\n# add/keep leading or trailing spaces).\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{ Old = \"\u062f\u0631\u0633 \u0627\u0645\u0631\u0648\u0632 \u062f\u0631\u0628\u0627\u0631\u0647 \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0632\u0645\u0627\u0646\u06cc \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f.\"; New = \"\u062f\u0631\u0633 \u0627\u0645\u0631\u0648\u0632 \u062f\u0631\u0628\u0627\u0631\u0647 \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0632\u0645\u0627\u0646\u06cc \u06a9\u0647 \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f.\" },\n    @{ Old = \"\u0627\u06cc\u0646\u062c\u0627 \u0633\u0647 \u062a\u0648\u0635\u06cc\u0647 \u0628\u0631\u0627\u06cc \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0648\u0642\u062a\u06cc \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f: \"; New = \"\u0627\u06cc\u0646\u062c\u0627 \u0633\u0647 \u062a\u0648\u0635\u06cc\u0647 \u0628\u0631\u0627\u06cc \u06a9\u0645\u06a9 \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a \u0648\u0642\u062a\u06cc \u06a9\u0647 \u0634\u062e\u0635\u06cc \u0648\u0641\u0627\u062a \u0645\u06cc\u06a9\u0646\u062f: \" },\n    @{ Old = \"\u0628\u0647 \u0632\u0628\u0627\u0646 \u0633\u0627\u062f\u0647\u060c \u0628\u0647 \u0622\u0646\u0647\u0627 \u0628\u06af\u0648\u06cc\u06cc\u062f \u06a9\u0647 \u0634\u062e\u0635 \u0648\u0641\u0627\u062a \u06a9\u0631\u062f\u0647 \u0648 \u062f\u06cc\u06af\u0631 \u0628\u0627\u0632 \u0646\u062e\u0648\u0627\u0647\u062f \u06af\u0634\u062a.\"; New = \"\u0628\u0647 \u0632\u0628\u0627\u0646 \u0633\u0627\u062f\u0647\u060c \u0628\u0647 \u0622\u0646\u0647\u0627 \u0628\u06af\u0648\u06cc\u06cc\u062f \u0634\u062e\u0635 \u06a9\u0647 \u0648\u0641\u0627\u062a \u06a9\u0631\u062f\u0647 \u0648 \u062f\u06cc\u06af\u0631 \u0628\u0627\u0632 \u0646\u062e\u0648\u0627\u0647\u062f \u06af\u0634\u062a.\" },\n    @{ Old = \"\u0645\u062a\u0648\u062c\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a \u0622\u0646\u200c \u0647\u0627 \u0628\u0627\u0634\u06cc\u062f \u0648 \u0622\u0646\u0686\u0647 \u0631\u0627 \u06a9\u0647 \u0645\u0634\u0627\u0647\u062f\u0647 \u0645\u06cc\u200c\u06a9\u0646\u06cc\u062f \u0628\u0647 \u0632\u0628\u0627\u0646 \u0628\u06cc\u0627\u0648\u0631\u06cc\u062f \u062a\u0627 \u0627\u062d\u0633\u0627\u0633 \u06a9\u0646\u0646\u062f \u06a9\u0647 \u06af\u067e \u0647\u0627\u06cc \u0634\u0627\u062a \u0634\u0646\u06cc\u062f\u0647 \u0645\u06cc\u200c\u0634\u0648\u062f.\"; New = \"\u0645\u062a\u0648\u062c\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a \u0622\u0646\u200c \u0647\u0627 \u0628\u0627\u0634\u06cc\u062f \u0648 \u0622\u0646\u0686\u0647 \u0631\u0627 \u06a9\u0647 \u0645\u0634\u0627\u0647\u062f\u0647 \u0645\u06cc\u200c\u06a9\u0646\u06cc\u062f \u0628\u0647 \u0632\u0628\u0627\u0646 \u0628\u06cc\u0627\u0648\u0631\u06cc\u062f \u062a\u0627 \u0627\u062d\u0633\u0627\u0633 \u06a9\u0646\u0646\u062f \u06a9\u0647 \u06af\u067e \u0647\u0627\u06cc \u0634\u0627\u0646 \u0634\u0646\u06cc\u062f\u0647 \u0645\u06cc\u200c\u0634\u0648\u062f.\" },\n    @{ Old = \"\u0645\u0645\u06a9\u0646 \u0627\u0633\u062a \u0644\u0627\u0632\u0645 \u0628\u0627\u0634\u062f \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u06a9\u0645\u06a9 \u06a9\u0646\u06cc\u062f \u062a\u0627 \u062f\u0631\u0628\u0627\u0631\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a\u0634 \u0635\u062d\u0628\u062a \u06a9\u0646\u0646\u062f \u0648 \u06cc\u0627 \u062a\u0648\u062c\u0647 \u0622\u0646\u0647\u0627 \u0631\u0627 \u0645\u062c\u062f\u062f\u0627 \u0645\u062a\u0645\u0631\u06a9\u0632 \u06a9\u0646\u06cc\u062f.\"; New = \"\u0645\u0645\u06a9\u0646 \u0627\u0633\u062a \u0644\u0627\u0632\u0645 \u0628\u0627\u0634\u062f \u0628\u0647 \u0627\u0637\u0641\u0627\u0644 \u06a9\u0645\u06a9 \u06a9\u0646\u06cc\u062f \u062a\u0627 \u062f\u0631\u0628\u0627\u0631\u0647 \u0627\u062d\u0633\u0627\u0633\u0627\u062a \u0634\u0627\u0646 \u0635\u062d\u0628\u062a \u06a9\u0646\u0646\u062f \u0648 \u06cc\u0627 \u062a\u0648\u062c\u0647 \u0622\u0646\u0647\u0627 \u0631\u0627 \u0645\u062c\u062f\u062f\u0627 \u0645\u062a\u0645\u0631\u06a9\u0632 \u06a9\u0646\u06cc\u062f.\" },\n    @{ Old = \"\u0634\u0634\u0645\u06cc\u0646 \u062a\u0648\u0635\u06cc\u0647 \u0622\u0631\u0627\u0645 \u06a9\u0631\u062f\u0646 \u0637\u0641\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a.\"; New = \"\u062a\u0648\u0635\u06cc\u0647 \u0634\u0634\u0645 \u0622\u0631\u0627\u0645 \u06a9\u0631\u062f\u0646 \u0637\u0641\u0644 \u062a\u0627\u0646 \u0627\u0633\u062a.\" }\n)\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    foreach ($edit in $edits) {\n        if ($r.Text.Contains($edit.Old)) {\n            $r.Text = $edit.New\n        }\n    }\n}\n\n$d.Save()\n"}
